$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.638.00'
$ws.Range("E2").Value = '  +4.15%  '
$ws.Range("D3").Value = '3.261.28'
$ws.Range("E3").Value = '  +4.39%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''579.53'
$ws.Range("E5").Value = '  +2.39%  '
$ws.Range("D6").Value = '''181.61'
$ws.Range("E6").Value = '  +7.93%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").Value = '''0.598'
$ws.Range("E8").Value = '  +1.75%  '
$ws.Range("D9").Value = '3.259.45'
$ws.Range("E9").Value = '  +4.37%  '
$ws.Range("E10").Value = '  +9.24%  '
$ws.Range("E11").Value = '  +3.33%  '
$ws.Range("E12").Value = '  +7.84%  '
$ws.Range("D13").Value = '3.832.47'
$ws.Range("E13").Value = '  +5.05%  '
$ws.Range("E14").Value = '  +1.57%  '
$ws.Range("D15").Value = '''28.39'
$ws.Range("E15").Value = '  +6.32%  '
$ws.Range("D16").Value = '67.624.30'
$ws.Range("E16").Value = '  +4.23%  '
$ws.Range("D17").Value = '''0.0000168'
$ws.Range("E17").Value = '  +4.98%  '
$ws.Range("D18").Value = '3.261.93'
$ws.Range("E18").Value = '  +4.57%  '
$ws.Range("D19").Value = '''5.83'
$ws.Range("E19").Value = '  +4.13%  '
$ws.Range("D20").Value = '''13.52'
$ws.Range("E20").Value = '  +7.46%  '
$ws.Range("D21").Value = '''375.01'
$ws.Range("E21").Value = '  +6.18%  '
$ws.Range("D22").Value = '''7.63'
$ws.Range("E22").Value = '  +6.83%  '
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("D24").Value = '''70.93'
$ws.Range("E24").Value = '  +3.70%  '
$ws.Range("D25").Value = '''0.510'
$ws.Range("E25").Value = '  +4.89%  '
$ws.Range("D26").Value = '''0.0000119'
$ws.Range("E26").Value = '  +7.77%  '
$ws.Range("D27").Value = '''9.60'
$ws.Range("E27").Value = '  +1.24%  '
$ws.Range("D28").Value = '''0.181'
$ws.Range("E28").Value = '  +3.69%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("E30").Value = '  +4.58%  '
$ws.Range("D31").Value = '''5.67'
$ws.Range("E31").Value = '  +9.00%  '
$ws.Range("D32").Value = '''22.71'
$ws.Range("E32").Value = '  +5.49%  '
$ws.Range("E33").Value = '  -0.07%  '
$ws.Range("E34").Value = '  +8.90%  '
$ws.Range("D35").Value = '''6.91'
$ws.Range("E35").Value = '  +6.52%  '
$ws.Range("D36").Value = '''163.55'
$ws.Range("E36").Value = '  +3.03%  '
$ws.Range("E37").Value = '  +6.98%  '
$ws.Range("D38").Value = '''0.848'
$ws.Range("E38").Value = '  +3.92%  '
$ws.Range("E39").Value = '  +6.55%  '
$ws.Range("D40").Value = '''6.83'
$ws.Range("E40").Value = '  +13.50%  '
$ws.Range("E41").Value = '  +3.93%  '
$ws.Range("E42").Value = '  +13.47%  '
$ws.Range("E43").Value = '  +8.35%  '
$ws.Range("D44").Value = '2.698.94'
$ws.Range("E44").Value = '  +3.10%  '
$ws.Range("D45").Value = '''351.96'
$ws.Range("E45").Value = '  +11.00%  '
$ws.Range("D46").Value = '''25.32'
$ws.Range("E46").Value = '  +8.46%  '
$ws.Range("D47").Value = '''40.83'
$ws.Range("E47").Value = '  +3.83%  '
$ws.Range("D48").Value = '''0.0679'
$ws.Range("E48").Value = '  +5.57%  '
$ws.Range("E49").Value = '  +4.39%  '
$ws.Range("D50").Value = '''0.998'
$ws.Range("E50").Value = '  +7.93%  '
$ws.Range("E51").Value = '  +1.63%  '
